$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values of the data rows (2-14), columns A-R (1-18),
# before overwriting anything, so the permutation below is computed
# from a consistent snapshot.
$colCount = 18
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14)

$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @()
    for ($c = 1; $c -le $colCount; $c++) {
        $rowVals += , $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Mapping: new row r <- old row mapping[r]
$mapping = @{
    2 = 5
    3 = 14
    4 = 12
    5 = 13
    6 = 8
    7 = 3
    8 = 6
    9 = 11
    10 = 4
    11 = 9
    12 = 2
    13 = 7
    14 = 10
}

foreach ($r in $rows) {
    $src = $snapshot[$mapping[$r]]
    for ($c = 1; $c -le $colCount; $c++) {
        $ws.Cells.Item($r, $c).Value = $src[$c - 1]
    }
}
